$d = $word.ActiveDocument

# Find the paragraph that starts the "# Modified productions" section.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*# Modified productions*") {
        $target = $i
        break
    }
}

# Delete everything from the very start of the document up to (but not
# including) that paragraph - i.e. the old "# Book productions" section
# plus the trailing blank paragraphs that followed it.
$startPara = $d.Paragraphs.Item(1)
$targetPara = $d.Paragraphs.Item($target)
$rangeToDelete = $d.Range($startPara.Range.Start, $targetPara.Range.Start)
$rangeToDelete.Delete()
